# The workbook gained one new data record. It was inserted as a new row
# at sheet row 307, which pushed the previously-existing rows 307-389
# down to rows 308-390 (dimension grows from A1:R389 to A1:R390).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 307; this shifts rows 307:389 -> 308:390
# and keeps their cell values/formatting completely intact.
$ws.Rows(307).Insert()

# Populate the newly inserted row 307 with its data.
$ws.Cells.Item(307, 1).Value  = 10
$ws.Cells.Item(307, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(307, 3).Value  = "La Araucanía"
$ws.Cells.Item(307, 4).Value  = 44642
$ws.Cells.Item(307, 5).Value  = 9
$ws.Cells.Item(307, 6).Value  = 100112032
$ws.Cells.Item(307, 7).Value  = "Zapallo italiano"
$ws.Cells.Item(307, 8).Value  = "Sin especificar"
$ws.Cells.Item(307, 9).Value  = "Primera"
$ws.Cells.Item(307, 10).Value = 125
$ws.Cells.Item(307, 11).Value = 13000
$ws.Cells.Item(307, 12).Value = 13000
$ws.Cells.Item(307, 13).Value = 13000
$ws.Cells.Item(307, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(307, 15).Value = "Región del Maule"
$ws.Cells.Item(307, 16).Value = 217
$ws.Cells.Item(307, 17).Value = 60
$ws.Cells.Item(307, 18).Value = "Hortaliza"
